# Applies the capital-structure database update described in the commit:
# refreshed financial metrics for existing companies, a company rename
# (Softlogic Capital -> Janashakthi Insurance) for the old row 5, and the
# insertion of a new Softlogic Capital PLC row, pushing Amana Takaful Life
# down to row 7. Sheet used-range grows from A1:AQ6 to A1:AQ7.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("AA2").Value = 0.05477415781032919
$ws.Range("AB2").Value = 0.1283252456543181
$ws.Range("AC2").Value = -0.07355108784398889
$ws.Range("AD2").Value = 47.609
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 47.609
$ws.Range("AG2").Value = 36.51000000000001
$ws.Range("AH2").Value = 0.2241699979753177
$ws.Range("AI2").Value = 0.1949917881380576
$ws.Range("AJ2").Value = 0.1813891096979332
$ws.Range("AK2").Value = 0.1566549386424097
$ws.Range("AL2").Value = 0.411
$ws.Range("AM2").Value = 0.411
$ws.Range("AN2").Value = 3.356291857596053
$ws.Range("AO2").Value = 33.2992700729927
$ws.Range("AP2").Value = 2.573845611561509
$ws.Range("AQ2").Value = 33.2992700729927
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "5"
$ws.Range("B2").Style = "Normal"
$ws.Range("D2").Value = 0.154
$ws.Range("E2").Value = -0.08539999999999999
$ws.Range("G2").Value = 0.06677648135438115
$ws.Range("H2").Value = 0.06677648135438115
$ws.Range("I2").Value = 0.06262182566918326
$ws.Range("J2").Value = 0.05344539068086505
$ws.Range("K2").Value = 7.164
$ws.Range("L2").Value = 0.03277968428277282
$ws.Range("M2").Value = 7.19
$ws.Range("N2").Value = 0.04363658432967166
$ws.Range("O2").Value = 1.003629257398102
$ws.Range("P2").Value = 7.19
$ws.Range("Q2").Value = 0.04363658432967166
$ws.Range("R2").Value = 1.003629257398102
$ws.Range("U2").Value = 11.099
$ws.Range("V2").Value = 0.06736056320932209
$ws.Range("W2").Value = 0.0551526717557252
$ws.Range("X2").Value = 0.1413494804821219
$ws.Range("Y2").Value = -0.08619680872639671
$ws.Range("Z2").Value = 1.156385917013239

# Row 3
$ws.Range("AA3").Value = 0.2319563127818913
$ws.Range("AB3").Value = 0.1268972576021954
$ws.Range("AC3").Value = 0.1050590551796959
$ws.Range("AD3").Value = 0.965
$ws.Range("AF3").Value = 0.965
$ws.Range("AG3").Value = -2.235
$ws.Range("AH3").Value = 0.1017395888244597
$ws.Range("AI3").Value = 0.09391727493917275
$ws.Range("AJ3").Value = -0.3556085918854416
$ws.Range("AK3").Value = -0.3159010600706714
$ws.Range("AL3").Value = 0.103
$ws.Range("AM3").Value = 0.103
$ws.Range("AN3").Value = 0.4974226804123711
$ws.Range("AO3").Value = 17.47572815533981
$ws.Range("AP3").Value = -1.152061855670103
$ws.Range("AQ3").Value = 17.47572815533981
$ws.Range("D3").Value = 0.314
$ws.Range("E3").Value = 0.391
$ws.Range("G3").Value = 0.181651376146789
$ws.Range("H3").Value = 0.181651376146789
$ws.Range("I3").Value = 0.1651376146788991
$ws.Range("J3").Value = 0.1185105234754452
$ws.Range("K3").Value = 1.22
$ws.Range("L3").Value = 0.1119266055045871
$ws.Range("U3").Value = 3.2
$ws.Range("V3").Value = 0.3755868544600939
$ws.Range("W3").Value = 0.148599269183922
$ws.Range("X3").Value = 0.1343014845048497
$ws.Range("Y3").Value = 0.01429778467907231
$ws.Range("Z3").Value = 1.957263422517507

# Row 4
$ws.Range("AA4").Value = 0.07997764662426574
$ws.Range("AB4").Value = 0.1250330368402513
$ws.Range("AC4").Value = -0.04505539021598556
$ws.Range("AD4").Value = 0.769
$ws.Range("AF4").Value = 0.769
$ws.Range("AG4").Value = -2.411
$ws.Range("AH4").Value = 0.007889688003365172
$ws.Range("AI4").Value = 0.009568365912229841
$ws.Range("AJ4").Value = -0.02557032103426699
$ws.Range("AK4").Value = -0.03123502053401392
$ws.Range("AN4").Value = 0.08277717976318623
$ws.Range("AP4").Value = -0.2595263724434876
$ws.Range("D4").Value = 0.128
$ws.Range("E4").Value = -0.135
$ws.Range("G4").Value = 0.1068616422947131
$ws.Range("H4").Value = 0.1068616422947131
$ws.Range("I4").Value = 0.09707536557930259
$ws.Range("J4").Value = 0.06607826934254576
$ws.Range("K4").Value = 6.31
$ws.Range("L4").Value = 0.07097862767154105
$ws.Range("M4").Value = 7.19
$ws.Range("N4").Value = 0.07435367114788004
$ws.Range("O4").Value = 1.13946117274168
$ws.Range("P4").Value = 7.19
$ws.Range("Q4").Value = 0.07435367114788004
$ws.Range("R4").Value = 1.13946117274168
$ws.Range("U4").Value = 3.18
$ws.Range("V4").Value = 0.03288521199586349
$ws.Range("W4").Value = 0.08538565629228687
$ws.Range("X4").Value = 0.1255380782952597
$ws.Range("Y4").Value = -0.04015242200297284
$ws.Range("Z4").Value = 1.210347174948945

# Row 5
$ws.Range("AA5").Value = 0.05477415781032919
$ws.Range("AB5").Value = 0.1283252456543181
$ws.Range("AC5").Value = -0.07355108784398889
$ws.Range("AD5").Value = 7.69
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 7.69
$ws.Range("AG5").Value = 6.510000000000001
$ws.Range("AH5").Value = 0.1736283585459472
$ws.Range("AI5").Value = 0.1275501741582352
$ws.Range("AJ5").Value = 0.151009046624913
$ws.Range("AK5").Value = 0.110133649128743
$ws.Range("AL5").Value = 0.263
$ws.Range("AM5").Value = 0.263
$ws.Range("AN5").Value = 2.496753246753247
$ws.Range("AO5").Value = 12.92775665399239
$ws.Range("AP5").Value = 2.113636363636364
$ws.Range("AQ5").Value = 12.92775665399239
$ws.Range("B5").Value = "Janashakthi Insurance PLC (COSE:JINS.N0000)"
$ws.Range("D5").Value = -0.09539999999999998
$ws.Range("E5").Value = -0.08539999999999999
$ws.Range("G5").Value = 0.09640718562874252
$ws.Range("H5").Value = 0.09640718562874252
$ws.Range("I5").Value = 0.1017964071856287
$ws.Range("J5").Value = 0.08845862491883702
$ws.Range("K5").Value = 2.89
$ws.Range("L5").Value = 0.08652694610778444
$ws.Range("U5").Value = 1.18
$ws.Range("V5").Value = 0.03224043715846994
$ws.Range("W5").Value = 0.0551526717557252
$ws.Range("X5").Value = 0.1423605509213477
$ws.Range("Y5").Value = -0.0872078791656225
$ws.Range("Z5").Value = 0.6192065257693734

# Row 6
$ws.Range("AA6").Value = 0
$ws.Range("AB6").Value = 0.1427711562217261
$ws.Range("AC6").Value = -0.1427711562217261
$ws.Range("AD6").Value = 37.7
$ws.Range("AF6").Value = 37.7
$ws.Range("AG6").Value = 34.39
$ws.Range("AH6").Value = 0.647766323024055
$ws.Range("AI6").Value = 0.417960088691796
$ws.Range("AJ6").Value = 0.6265257788303881
$ws.Range("AK6").Value = 0.3957877776499022
$ws.Range("AL6").Value = 0
$ws.Range("AM6").Value = 0
$ws.Range("B6").Value = "Softlogic Capital PLC (COSE:SCAP.N0000)"
$ws.Range("D6").Value = 0.18
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = -3.08
$ws.Range("L6").Value = -0.03869346733668342
$ws.Range("U6").Value = 3.31
$ws.Range("V6").Value = 0.1614634146341463
$ws.Range("W6").Value = -0.09967637540453075
$ws.Range("X6").Value = 0.2779106369732995
$ws.Range("Y6").Value = -0.3775870123778302
$ws.Range("Z6").Value = 1.494274450910456

# Row 7
$ws.Range("A7").Value = "Sri Lanka"
$ws.Range("AA7").Value = -0.05207956600361662
$ws.Range("AB7").Value = 0.1289440137740778
$ws.Range("AC7").Value = -0.1810235797776944
$ws.Range("AD7").Value = 0.485
$ws.Range("AE7").Value = 0
$ws.Range("AF7").Value = 0.485
$ws.Range("AG7").Value = 0.256
$ws.Range("AH7").Value = 0.1652470187393526
$ws.Range("AI7").Value = 0.1603305785123967
$ws.Range("AJ7").Value = 0.09460458240946044
$ws.Range("AK7").Value = 0.09155937052932761
$ws.Range("AL7").Value = 0.045
$ws.Range("AM7").Value = 0.045
$ws.Range("AN7").Value = -3.88
$ws.Range("AO7").Value = -3.2
$ws.Range("AP7").Value = -2.048
$ws.Range("AQ7").Value = -3.2
$ws.Range("B7").Value = "Amãna Takaful Life PLC (COSE:ATLL.N0000)"
$ws.Range("C7").Value = "Insurance (Life)"
$ws.Range("G7").Value = -0.01843478260869565
$ws.Range("H7").Value = -0.01843478260869565
$ws.Range("I7").Value = -0.02504347826086956
$ws.Range("J7").Value = -0.02504347826086956
$ws.Range("K7").Value = -0.176
$ws.Range("L7").Value = -0.03060869565217391
$ws.Range("M7").Value = -0
$ws.Range("N7").Value = -0
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = -0
$ws.Range("Q7").Value = -0
$ws.Range("R7").Value = 0
$ws.Range("S7").Value = 0
$ws.Range("U7").Value = 0.229
$ws.Range("V7").Value = 0.09346938775510204
$ws.Range("W7").Value = -0.06423357664233575
$ws.Range("X7").Value = 0.1413494804821219
$ws.Range("Y7").Value = -0.2055830571244577
$ws.Range("Z7").Value = 2.079566003616636

# Row 6 previously ended at column AM; the new Softlogic Capital PLC row
# has no data beyond AM, so clear the stale AN:AQ values left over from
# the row that used to occupy row 6 (Amana Takaful Life, now row 7).
$ws.Range("AN6").ClearContents()
$ws.Range("AO6").ClearContents()
$ws.Range("AP6").ClearContents()
$ws.Range("AQ6").ClearContents()

Write-Output "done"
